# Apply scheduled market-data refresh updates to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3138.6365
$ws.Range("I64").Value = 2805
$ws.Range("K64").Value = 2805
$ws.Range("M64").Value = -2557

$ws.Range("H67").Value = 3138.6365
$ws.Range("I67").Value = 2805
$ws.Range("K67").Value = 2805
$ws.Range("M67").Value = -1947

$ws.Range("H74").Value = 6254366.5
$ws.Range("J74").Value = 11368216
$ws.Range("L74").Value = 11368216
$ws.Range("N74").Value = -11370088

$ws.Range("H77").Value = 6254366.5
$ws.Range("J77").Value = 11368216
$ws.Range("L77").Value = 56841080
$ws.Range("N77").Value = -56850440

$ws.Range("H111").Value = 2526.8
$ws.Range("I111").Value = 4819.5
$ws.Range("J111").Value = 998.3333
$ws.Range("K111").Value = 14458.5
$ws.Range("L111").Value = 2994.9999
$ws.Range("M111").Value = -11391.5
$ws.Range("N111").Value = -9128.999899999999

$ws.Range("H129").Value = 264355.1
$ws.Range("I129").Value = 296.85715
$ws.Range("J129").Value = 323981.16
$ws.Range("K129").Value = 890.5714499999999
$ws.Range("L129").Value = 971943.48
$ws.Range("M129").Value = 4109.428550000001
$ws.Range("N129").Value = -981943.48

$ws.Range("H132").Value = 3039.139
$ws.Range("I132").Value = 3116.9666
$ws.Range("J132").Value = 2650
$ws.Range("K132").Value = 9350.899800000001
$ws.Range("L132").Value = 7950
$ws.Range("M132").Value = -6820.899800000001
$ws.Range("N132").Value = -13010

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 692.9211
$ws.Range("I2").Value = 615.4483
$ws.Range("K2").Value = 615.4483
$ws.Range("M2").Value = -502.4483

$ws.Range("H32").Value = 5846.6704
$ws.Range("I32").Value = 4505.896
$ws.Range("J32").Value = 18751.625
$ws.Range("K32").Value = 4505.896
$ws.Range("L32").Value = 18751.625
$ws.Range("M32").Value = -4218.896
$ws.Range("N32").Value = -19325.625

$ws.Range("H45").Value = 2418
$ws.Range("I45").Value = 2272.5
$ws.Range("J45").Value = 2650.8
$ws.Range("K45").Value = 2272.5
$ws.Range("L45").Value = 2650.8
$ws.Range("M45").Value = -1895.5
$ws.Range("N45").Value = -3404.8

$ws.Range("H116").Value = 692.9211
$ws.Range("I116").Value = 615.4483
$ws.Range("K116").Value = 615.4483
$ws.Range("M116").Value = 1678.5517

$ws.Range("H122").Value = 2239.2856
$ws.Range("I122").Value = 2043.421
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 6130.263
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -3680.263
$ws.Range("N122").Value = -17200

$ws.Range("H132").Value = 10779.105
$ws.Range("I132").Value = 2004.525
$ws.Range("J132").Value = 31425.176
$ws.Range("K132").Value = 6013.575000000001
$ws.Range("L132").Value = 94275.52799999999
$ws.Range("M132").Value = -3483.575000000001
$ws.Range("N132").Value = -99335.52799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 692.9211
$ws.Range("I3").Value = 615.4483
$ws.Range("K3").Value = 615.4483
$ws.Range("M3").Value = -501.4483

$ws.Range("H82").Value = 35889.4
$ws.Range("J82").Value = 42797.5
$ws.Range("L82").Value = 42797.5
$ws.Range("N82").Value = -43563.5

$ws.Range("H85").Value = 35889.4
$ws.Range("J85").Value = 42797.5
$ws.Range("L85").Value = 42797.5
$ws.Range("N85").Value = -45449.5

$ws.Range("H86").Value = 1881.2162
$ws.Range("I86").Value = 1765.4348
$ws.Range("J86").Value = 2071.4285
$ws.Range("K86").Value = 1765.4348
$ws.Range("L86").Value = 2071.4285
$ws.Range("M86").Value = -642.4348
$ws.Range("N86").Value = -4317.4285

$ws.Range("H89").Value = 1881.2162
$ws.Range("I89").Value = 1765.4348
$ws.Range("J89").Value = 2071.4285
$ws.Range("K89").Value = 8827.173999999999
$ws.Range("L89").Value = 10357.1425
$ws.Range("M89").Value = -3211.173999999999
$ws.Range("N89").Value = -21589.1425

$ws.Range("H134").Value = 3171.8865
$ws.Range("I134").Value = 3426.4595
$ws.Range("J134").Value = 1826.2858
$ws.Range("K134").Value = 10279.3785
$ws.Range("L134").Value = 5478.857400000001
$ws.Range("M134").Value = -7744.378499999999
$ws.Range("N134").Value = -10548.8574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 984.36365
$ws.Range("I16").Value = 872
$ws.Range("J16").Value = 1181
$ws.Range("K16").Value = 872
$ws.Range("L16").Value = 1181
$ws.Range("M16").Value = -585
$ws.Range("N16").Value = -1755

$ws.Range("H31").Value = 3789.0852
$ws.Range("I31").Value = 2330.0588
$ws.Range("J31").Value = 4615.8667
$ws.Range("K31").Value = 2330.0588
$ws.Range("L31").Value = 4615.8667
$ws.Range("M31").Value = -2035.0588
$ws.Range("N31").Value = -5205.8667

$ws.Range("H34").Value = 3789.0852
$ws.Range("I34").Value = 2330.0588
$ws.Range("J34").Value = 4615.8667
$ws.Range("K34").Value = 2330.0588
$ws.Range("L34").Value = 4615.8667
$ws.Range("M34").Value = -2128.0588
$ws.Range("N34").Value = -5019.8667

$ws.Range("H99").Value = 3337.6296
$ws.Range("I99").Value = 2600.7273
$ws.Range("J99").Value = 6580
$ws.Range("K99").Value = 2600.7273
$ws.Range("L99").Value = 6580
$ws.Range("M99").Value = -1102.7273
$ws.Range("N99").Value = -9576

$ws.Range("H113").Value = 984.36365
$ws.Range("I113").Value = 872
$ws.Range("J113").Value = 1181
$ws.Range("K113").Value = 872
$ws.Range("L113").Value = 1181
$ws.Range("M113").Value = 1298
$ws.Range("N113").Value = -5521

$ws.Range("H122").Value = 986.7143
$ws.Range("I122").Value = 794.2727
$ws.Range("J122").Value = 1198.4
$ws.Range("K122").Value = 2382.8181
$ws.Range("L122").Value = 3595.2
$ws.Range("M122").Value = 67.18190000000004
$ws.Range("N122").Value = -8495.200000000001

$ws.Range("H126").Value = 3337.6296
$ws.Range("I126").Value = 2600.7273
$ws.Range("J126").Value = 6580
$ws.Range("K126").Value = 7802.1819
$ws.Range("L126").Value = 19740
$ws.Range("M126").Value = -5332.1819
$ws.Range("N126").Value = -24680

$ws.Range("H134").Value = 665.65955
$ws.Range("I134").Value = 591.7632
$ws.Range("K134").Value = 1775.2896
$ws.Range("M134").Value = 759.7103999999999

$ws.Range("H141").Value = 25112.79
$ws.Range("J141").Value = 25112.79
$ws.Range("L141").Value = 25112.79
$ws.Range("N141").Value = -35472.79

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 39
$ws.Range("I61").Value = 39
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 117
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = 98

$ws.Range("H114").Value = 2799.75
$ws.Range("J114").Value = 999.5
$ws.Range("L114").Value = 2998.5
$ws.Range("N114").Value = -9506.5

$ws.Range("H131").Value = 716.83
$ws.Range("J131").Value = 739.2526
$ws.Range("L131").Value = 2217.7578
$ws.Range("N131").Value = -12297.7578

$ws.Range("H132").Value = 912.7917
$ws.Range("I132").Value = 740.8
$ws.Range("J132").Value = 1035.6428
$ws.Range("K132").Value = 6667.2
$ws.Range("L132").Value = 9320.7852
$ws.Range("M132").Value = -4137.2
$ws.Range("N132").Value = -14380.7852

$ws.Range("H140").Value = 3145.8572
$ws.Range("I140").Value = 1534
$ws.Range("J140").Value = 5295
$ws.Range("K140").Value = 4602
$ws.Range("L140").Value = 15885
$ws.Range("M140").Value = 578
$ws.Range("N140").Value = -26245

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1256.6
$ws.Range("I97").Value = 1336
$ws.Range("J97").Value = 995.7143
$ws.Range("K97").Value = 1336
$ws.Range("L97").Value = 995.7143
$ws.Range("M97").Value = -840
$ws.Range("N97").Value = -1987.7143

$ws.Range("H102").Value = 2283.25
$ws.Range("I102").Value = 1759.5
$ws.Range("K102").Value = 1759.5
$ws.Range("M102").Value = -137.5

$ws.Range("H113").Value = 6221.1055
$ws.Range("I113").Value = 8576.833000000001
$ws.Range("J113").Value = 2182.7144
$ws.Range("K113").Value = 8576.833000000001
$ws.Range("L113").Value = 2182.7144
$ws.Range("M113").Value = -6406.833000000001
$ws.Range("N113").Value = -6522.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1631.5294
$ws.Range("I93").Value = 1445.4286
$ws.Range("K93").Value = 1445.4286
$ws.Range("M93").Value = -197.4286

$ws.Range("H100").Value = 1990.8
$ws.Range("I100").Value = 879.3333
$ws.Range("J100").Value = 2731.7778
$ws.Range("K100").Value = 879.3333
$ws.Range("L100").Value = 2731.7778
$ws.Range("M100").Value = -338.3333
$ws.Range("N100").Value = -3813.7778

$ws.Range("H122").Value = 615777.4
$ws.Range("I122").Value = 703053.6
$ws.Range("J122").Value = 4843.75
$ws.Range("K122").Value = 2109160.8
$ws.Range("L122").Value = 14531.25
$ws.Range("M122").Value = -2106710.8
$ws.Range("N122").Value = -19431.25

$ws.Range("H130").Value = 19800
$ws.Range("J130").Value = 19800
$ws.Range("L130").Value = 19800
$ws.Range("N130").Value = -29840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 248.14285
$ws.Range("I100").Value = 256.72726
$ws.Range("K100").Value = 513.45452
$ws.Range("M100").Value = 27.54548

$ws.Range("H126").Value = 1843.0294
$ws.Range("I126").Value = 1399
$ws.Range("J126").Value = 3286.125
$ws.Range("K126").Value = 4197
$ws.Range("L126").Value = 9858.375
$ws.Range("M126").Value = -1727
$ws.Range("N126").Value = -14798.375
